# Apply updates to FAST_holdings workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow edits, will re-protect at the end
$ws.Unprotect()

# Update the confidentiality / "as of" date string directly in its known cell (A13)
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."
$ws.Range("A13").Value = $newText
# Restore auto row height (setting the multi-line text can otherwise leave an explicit row height)
$ws.Rows("13:13").AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.09337966529029698
$ws.Range("E2").Value = 0.006185822095756466

$ws.Range("D3").Value = 0.1075231391539414
$ws.Range("E3").Value = 0.005606704438149146

$ws.Range("D4").Value = 0.1193736659991157
$ws.Range("E4").Value = 0.0004661642452024051

$ws.Range("D5").Value = 0.1402420547274613
$ws.Range("E5").Value = 0.005680473372780881

$ws.Range("D6").Value = 0.135399725514456
$ws.Range("E6").Value = -0.0005527915975677145

$ws.Range("D7").Value = 0.1455764270701752
$ws.Range("E7").Value = 0.009341845836071672

$ws.Range("D8").Value = 0.1279149424567194
$ws.Range("E8").Value = 0.01111761263897026

$ws.Range("D9").Value = 0.1305903797878341
$ws.Range("E9").Value = 0.012883120793868

$ws.Range("E10").Value = 0.00642239457764826

# Restore sheet protection (matching original intent: sheet/objects/scenarios
# protected, but row/column formatting still allowed)
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true)

$wb.Save()
